$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.428.11"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "1.916.75"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9978"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4645"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4153"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08055"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.023"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").Value = "1.883.13"
$ws.Range("E13").Value = "  -3.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.993"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.169"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9979"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001034"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06589"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "29.370.14"
$ws.Range("E22").Value = "  -1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.539"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "
$ws.Range("D26").Value = "2.128.61"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.163"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.665"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.043"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09461"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.438"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.446"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.531"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06131"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02267"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.463"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5915"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9977"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1839"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.374"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.74%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.243"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07531"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("B49").Value = "Decentraland"
$ws.Range("C49").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5581"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.936"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.09%  "